# Change the East-Asian font used throughout the style sheet from
# "DejaVu Sans" to "Tahoma", and make sure the complex-script ("cs")
# font is explicitly recorded as "DejaVu Sans" on the paragraph styles
# that currently inherit it implicitly (List, Caption, Index).

$d = $word.ActiveDocument

# docDefaults (w:rPrDefault) and the "Normal" style share the same
# East-Asian font value in this template; update the visible style
# that the object model exposes (the document default inherits from
# it for any run that does not specify its own East-Asian font).
$normal = $d.Styles("Normal")
$normal.Font.NameFarEast = "Tahoma"

$heading = $d.Styles("Heading")
$heading.Font.NameFarEast = "Tahoma"

# These styles had an empty <w:rPr/> (or, for Caption, an <w:rPr> with
# no rFonts element) and picked up their complex-script font from the
# style they are based on. Setting the bidi/complex-script font name
# explicitly materializes a <w:rFonts w:cs="DejaVu Sans"/> element.
$list = $d.Styles("List")
$list.Font.NameBi = "DejaVu Sans"

$caption = $d.Styles("Caption")
$caption.Font.NameBi = "DejaVu Sans"

$index = $d.Styles("Index")
$index.Font.NameBi = "DejaVu Sans"

Write-Output "Styles updated."
